$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 was "Bowel" / "Coming soon" placeholder; fill in the new entry.
$ws.Range("B2").Value = "Appendiceal Mucocele"
$ws.Range("C2").Value = "Clip 1 B-mode + Color"
$ws.Range("D2").Value = "https://youtu.be/kdZO1IPuOIw"

# Turn D2 into a real hyperlink pointing at the YouTube clip, then restore
# the standard "hyperlink" cell style used by the rest of column D.
$ws.Hyperlinks.Add($ws.Range("D2"), "https://youtu.be/kdZO1IPuOIw", "", "", "https://youtu.be/kdZO1IPuOIw")
$ws.Range("D2").Style = "Collegamento ipertestuale"

# Leave the selection where the author left it when they saved.
[void]$ws.Range("D7").Select()
